# ---------------------------------------------------------------------------
# Applies the recorded edit to the deck:
#   1. Slide 5's table switches from the custom "Table_0" style to the
#      built-in table style {17C84CA8-C28F-4469-B706-E559B8F11C29}.
#   2. The presentation's theme (used by the slide master / all slides)
#      is re-coloured from the "Red Violet" / "Integral" scheme to the
#      standard Office "Office Theme" colour scheme.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------------
$slide  = $p.Slides.Item(5)
$tShape = $slide.Shapes.Item(2)
$table  = $tShape.Table
$table.ApplyStyle("{17C84CA8-C28F-4469-B706-E559B8F11C29}")

# --- 2. Theme colour scheme -------------------------------------------------
function ToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$colors.Item(1).RGB  = ToRgbLong("000000")
$colors.Item(2).RGB  = ToRgbLong("FFFFFF")
$colors.Item(3).RGB  = ToRgbLong("44546A")
$colors.Item(4).RGB  = ToRgbLong("E7E6E6")
$colors.Item(5).RGB  = ToRgbLong("5B9BD5")
$colors.Item(6).RGB  = ToRgbLong("ED7D31")
$colors.Item(7).RGB  = ToRgbLong("A5A5A5")
$colors.Item(8).RGB  = ToRgbLong("FFC000")
$colors.Item(9).RGB  = ToRgbLong("4472C4")
$colors.Item(10).RGB = ToRgbLong("70AD47")
$colors.Item(11).RGB = ToRgbLong("0563C1")
$colors.Item(12).RGB = ToRgbLong("954F72")

Write-Output "table style + theme colours updated"
